# correct size of scratch-nompiio
#
# The "scratch-nompiio" capacity label on the storage-overview slide read
# "250 Tb    /work/scratch-nompiio" but the real size is 175 Tb. Locate the
# textbox (wherever it lives) and fix just the leading number, which splits
# the original single run into two runs ("175 " / "Tb    /work/scratch-")
# exactly like a user retyping the number in-place would.

$p = $ppt.ActivePresentation

$targetOld = "250 Tb    /work/scratch-nompiio"
$findWhat  = "250 "
$replaceWhat = "175 "

$fixed = $false

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        if (-not $shape.HasTextFrame) { continue }
        if (-not $shape.TextFrame.HasText) { continue }

        $tr = $shape.TextFrame.TextRange

        if ($tr.Text -eq $targetOld) {
            [void]$tr.Replace($findWhat, $replaceWhat, 0, 0, 0)
            $fixed = $true
            break
        }
    }

    if ($fixed) { break }
}
